$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '57.813.02'
$ws.Range("E2").Value = '  +0.64%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.946.08'
$ws.Range("E3").Value = '  +2.05%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '553.01'
$ws.Range("E5").Value = '  +0.50%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '133.26'
$ws.Range("E6").Value = '  +9.74%  '
$ws.Range("E7").Value = '  -0.12%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.514'
$ws.Range("E8").Value = '  +5.26%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.940.88'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.128'
$ws.Range("E10").Value = '  +3.24%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '4.80'
$ws.Range("E11").Value = '  +0.32%  '
$ws.Range("E12").Value = '  +4.48%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000220'
$ws.Range("E13").Value = '  +4.87%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '32.81'
$ws.Range("E14").Value = '  +4.98%  '
$ws.Range("E15").Value = '  +2.96%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.432.95'
$ws.Range("E16").Value = '  +2.14%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.94'
$ws.Range("E17").Value = '  +9.47%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.943.94'
$ws.Range("E18").Value = '  +2.18%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '57.806.69'
$ws.Range("E19").Value = '  +0.77%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '415.62'
$ws.Range("E20").Value = '  +2.36%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.36'
$ws.Range("E21").Value = '  +5.19%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.697'
$ws.Range("E22").Value = '  +7.53%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '13.44'
$ws.Range("E23").Value = '  +7.42%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.01'
$ws.Range("E24").Value = '  +4.13%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '78.98'
$ws.Range("E25").Value = '  +3.29%  '
$ws.Range("E27").Value = '  +0.20%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.49'
$ws.Range("E28").Value = '  +1.34%  '
$ws.Range("E29").Value = '  +6.66%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.49'
$ws.Range("E30").Value = '  +5.94%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '25.45'
$ws.Range("E31").Value = '  +3.88%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.95'
$ws.Range("E32").Value = '  -0.69%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0963'
$ws.Range("E33").Value = '  +1.85%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.67'
$ws.Range("E34").Value = '  +6.30%  '
$ws.Range("E35").Value = '  +5.90%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.06'
$ws.Range("E36").Value = '  +2.62%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '8.99'
$ws.Range("E37").Value = '  +7.89%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0₃0696'
$ws.Range("E38").Value = '  +13.87%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '48.34'
$ws.Range("E39").Value = '  +0.00%  '
$ws.Range("E40").Value = '  +15.75%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '381.43'
$ws.Range("E42").Value = '  +2.96%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0347'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.706.92'
$ws.Range("E44").Value = '  +3.85%  '
$ws.Range("E45").Value = '  +0.03%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '124.42'
$ws.Range("E46").Value = '  +5.84%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.237'
$ws.Range("E47").Value = '  +4.29%  '
$ws.Range("E48").Value = '  +2.64%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.97'
$ws.Range("E49").Value = '  +2.62%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '22.88'
$ws.Range("E50").Value = '  +2.82%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.99'
$ws.Range("E51").Value = '  +3.53%  '
